$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of E4 and E5 (circleCount now controlled by rotation instead of Mask)
$ws.Range("E4").Value = 15
$ws.Range("E5").Value = 1

# Update the active selection to E4
$ws.Range("E4").Select()
